# Generate Report for Handoff
#
# The localization CI regenerated the status report: the "In Translation"
# status became "Ready for handoff" (its timestamp bumped a few seconds
# later) everywhere it appears, the zh-cn handoff timestamp moved forward,
# and the columns holding the longer status text were re-autosized.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" everywhere it shows up
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps
# Overview's "Latest HO Xliff Generate Date" and de-de's "Latest Handoff
# Datetime" share the same value, bumped by 34s.
$overview.Range("G2").Value = "2016-08-24 04:56:08"
$dede.Range("H2").Value     = "2016-08-24 04:56:08"

# zh-cn's "Latest Handoff Datetime" bumped by 29s.
$zhcn.Range("H2").Value = "2016-08-24 04:55:59"

# --- Column re-sizing to fit the new, longer "Ready for handoff" text
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth     = 16.33
$dede.Columns.Item(3).ColumnWidth     = 16.33
